$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matches")

$homeFormula = "=INDEX(seeds[team],MATCH(matches[[#This Row],[home-seed]],seeds[seed],0))"
$awayFormula = "=INDEX(seeds[team],MATCH(matches[[#This Row],[away-seed]],seeds[seed],0))"

for ($r = 2; $r -le 25; $r++) {
  $ws.Cells.Item($r, 6).Formula = $homeFormula
  $ws.Cells.Item($r, 7).Formula = $awayFormula
}

$ws.Activate()
$ws.Range("F26:G33").Select()
